$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEAM")

$ws.Range("A3").Value = "Ni"
$ws.Range("B3").Value = 4.45
$ws.Range("C3").Value = 2.49
$ws.Range("D3").Formula = "=F5"
$ws.Range("E3").Value = 0.94
$ws.Range("F3").Value = 2.56
$ws.Range("G3").Value = 1.5
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 1.5
$ws.Range("K3").Value = 3.1
$ws.Range("L3").Value = 1.8
$ws.Range("M3").Value = 4.36
$ws.Range("N3").Value = 0.81
$ws.Range("R3").Value = 0.05
$ws.Range("S3").Value = 0.05
$ws.Range("T3").Value = 1

$ws.Range("C7").Value = "fcc"
$ws.Range("D5").Formula = "=0.0062415*187.6"
$ws.Range("S7").ClearContents()
